$d = $word.ActiveDocument

# --- 1. Paragraph 36: "If you had 400 points..." -> single run, lang en-CA, drop inline bookmark ---
$p36 = $d.Paragraphs(36)
$xmlP36 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DC5CF5" w:rsidRPr="00D57C17" w:rsidRDefault="00D57C17" w:rsidP="00D57C17"><w:pPr><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr><w:t>If you had 400 points to distribute between the members of your team depending on their participation and implication, how would you allocate them?</w:t></w:r></w:p>'
$p36.Range.InsertXML($xmlP36)

# --- 2. Paragraph 34: Paragraphedeliste empty paragraph -> lang en-CA ---
$p34 = $d.Paragraphs(34)
$xmlP34 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D57C17" w:rsidRPr="00D57C17" w:rsidRDefault="00D57C17" w:rsidP="00D57C17"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr></w:pPr></w:p>'
$p34.Range.InsertXML($xmlP34)

# --- 3. Paragraph 33: "What is your opinion..." -> single run, lang en-CA ---
$p33 = $d.Paragraphs(33)
$xmlP33 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D57C17" w:rsidRPr="00D57C17" w:rsidRDefault="00D57C17" w:rsidP="00D57C17"><w:pPr><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr><w:t>What is your opinion on the algorithmic, maven, git, demo sessions? Did it help you in any way? (half a page)</w:t></w:r></w:p>'
$p33.Range.InsertXML($xmlP33)

# --- 4. Paragraph 29: drop <w:lastRenderedPageBreak/> ---
$p29 = $d.Paragraphs(29)
$xmlP29 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="006E4F08" w:rsidRDefault="006E4F08" w:rsidP="006E4F08"><w:pPr><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr></w:pPr><w:r w:rsidRPr="006E4F08"><w:rPr><w:lang w:val="en-CA" w:eastAsia="fr-FR"/></w:rPr><w:t xml:space="preserve">How have you been able to benefit from the milestones you defined? How would you modify them, now that you are at the end of the project? How have you integrated the feedback given every day? </w:t></w:r></w:p>'
$p29.Range.InsertXML($xmlP29)

# --- 5. Paragraph 7: delete the "Le " run paragraph, merging it into the following blank one ---
$p7 = $d.Paragraphs(7)
$p7.Range.Delete()

# --- 6. The now-blank paragraph (was paragraph 8) gets the relocated _GoBack bookmark ---
$pBlank = $d.Paragraphs(7)
$d.Bookmarks.Add("_GoBack", $pBlank.Range)

# --- 7. Paragraph 6: "Problème n-p complet." -> expanded sentence with spell-checked "Hashcode" ---
$p6 = $d.Paragraphs(6)
$xmlP6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00413229" w:rsidRDefault="00413229" w:rsidP="0076465F"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr><w:t xml:space="preserve">Le problème énoncé par le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr><w:t>Hashcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr><w:t xml:space="preserve"> de Google est un p</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr><w:t>roblème NP</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="fr-FR"/></w:rPr><w:t xml:space="preserve"> complet.</w:t></w:r></w:p>'
$p6.Range.InsertXML($xmlP6)
